# Frame in main application with self definition
#
# - Adds a new "Startup" worksheet after Sheet2, populated with the p1-p8 /
#   v1-v8 parameter labels (column A) and zeroed values (column B).
# - Makes "Startup" the active/selected sheet.
# - Updates the selection remembered on Sheet2 (was F12, now G29).

$wb = $excel.ActiveWorkbook

$ws2 = $wb.Worksheets.Item("Sheet2")

# Sheet2 now remembers a different selected cell.
$ws2.Activate()
[void]$ws2.Range("G29").Select()

# Insert the new "Startup" sheet right after Sheet2.
$startup = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws2)
$startup.Name = "Startup"

$labels = @("p1", "p2", "p3", "p4", "p5", "p6", "p7", "p8", "v1", "v2", "v3", "v4", "v5", "v6", "v7", "v8")

for ($i = 0; $i -lt $labels.Length; $i++) {
    $row = $i + 1
    $startup.Cells.Item($row, 1).Value = $labels[$i]
    $startup.Cells.Item($row, 2).Value = 0
}

# Leave the new sheet active, with the same lingering selection as the source file.
[void]$startup.Range("K8").Select()
